# Clear the "Mês de Competencia" column (F) header and its sample value,
# leaving the cell styles intact but with no content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = ""
$ws.Range("F2").Value = ""

# Update the active selection to match the saved workbook state.
$ws.Range("F5").Select()
